$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 25, pushing the existing rows 25-27 down to 28-30.
$ws.Rows("25:27").Insert()

# Row 25 (new): same as the old row 25 (now row 28) but with an updated
# date, quality ("Calidad") and volume.
$ws.Cells.Item(25, 1).Value = 3
$ws.Cells.Item(25, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 44455
$ws.Cells.Item(25, 5).Value = 5
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100107
$ws.Cells.Item(25, 8).Value = "Otros"
$ws.Cells.Item(25, 9).Value = 100107002
$ws.Cells.Item(25, 10).Value = "Chirimoya"
$ws.Cells.Item(25, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 45
$ws.Cells.Item(25, 14).Value = 30000
$ws.Cells.Item(25, 15).Value = 30000
$ws.Cells.Item(25, 16).Value = 30000
$ws.Cells.Item(25, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(25, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(25, 19).Value = 3000
$ws.Cells.Item(25, 20).Value = 10

# Row 26 (new): same as the old row 26 (now row 29) but with an updated
# date, volume and resulting prices.
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44455
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100107
$ws.Cells.Item(26, 8).Value = "Otros"
$ws.Cells.Item(26, 9).Value = 100107002
$ws.Cells.Item(26, 10).Value = "Chirimoya"
$ws.Cells.Item(26, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 48
$ws.Cells.Item(26, 14).Value = 27000
$ws.Cells.Item(26, 15).Value = 27000
$ws.Cells.Item(26, 16).Value = 27000
$ws.Cells.Item(26, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(26, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(26, 19).Value = 2700
$ws.Cells.Item(26, 20).Value = 10

# Row 27 (new): brand-new record.
$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(27, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = 44455
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100107
$ws.Cells.Item(27, 8).Value = "Otros"
$ws.Cells.Item(27, 9).Value = 100107002
$ws.Cells.Item(27, 10).Value = "Chirimoya"
$ws.Cells.Item(27, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 12).Value = "Segunda"
$ws.Cells.Item(27, 13).Value = 47
$ws.Cells.Item(27, 14).Value = 25000
$ws.Cells.Item(27, 15).Value = 25000
$ws.Cells.Item(27, 16).Value = 25000
$ws.Cells.Item(27, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(27, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(27, 19).Value = 2500
$ws.Cells.Item(27, 20).Value = 10
